$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2182.9546
$ws.Range("J2").Value = 3381.818
$ws.Range("L2").Value = 3381.818
$ws.Range("N2").Value = -3607.818
$ws.Range("H9").Value = 90.62069
$ws.Range("I9").Value = 74.652176
$ws.Range("J9").Value = 151.83333
$ws.Range("K9").Value = 74.652176
$ws.Range("L9").Value = 151.83333
$ws.Range("M9").Value = 94.347824
$ws.Range("N9").Value = -489.83333
$ws.Range("H15").Value = 230.59
$ws.Range("I15").Value = 230.59
$ws.Range("K15").Value = 691.77
$ws.Range("M15").Value = -522.77
$ws.Range("H52").Value = 2333.3333
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 3000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("M52").Value = -2840
$ws.Range("N52").Value = -9320
$ws.Range("H118").Value = 1840
$ws.Range("J118").Value = 1867.5883
$ws.Range("L118").Value = 5602.7649
$ws.Range("N118").Value = -8916.7649
$ws.Range("H129").Value = 325462.53
$ws.Range("J129").Value = 670283
$ws.Range("L129").Value = 2010849
$ws.Range("N129").Value = -2020849
$ws.Range("H135").Value = 31251996
$ws.Range("I135").Value = 1090.6364
$ws.Range("J135").Value = 100003980
$ws.Range("K135").Value = 9815.7276
$ws.Range("L135").Value = 900035820
$ws.Range("M135").Value = -7280.7276
$ws.Range("N135").Value = -900040890
$ws.Range("H137").Value = 32631056
$ws.Range("I137").Value = 12501060
$ws.Range("K137").Value = 37503180
$ws.Range("M137").Value = -37500630
$ws.Range("H138").Value = 3266.9155
$ws.Range("I138").Value = 3285.842
$ws.Range("J138").Value = 3260
$ws.Range("K138").Value = 9857.526
$ws.Range("L138").Value = 9780
$ws.Range("M138").Value = -4717.526
$ws.Range("N138").Value = -20060

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 100007
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 100007
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100007
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -100525
$ws.Range("H37").Value = 9857.6
$ws.Range("J37").Value = 13468.363
$ws.Range("L37").Value = 13468.363
$ws.Range("N37").Value = -14014.363
$ws.Range("H44").Value = 22487
$ws.Range("J44").Value = 22487
$ws.Range("L44").Value = 22487
$ws.Range("N44").Value = -23463
$ws.Range("H55").Value = 24106
$ws.Range("J55").Value = 24106
$ws.Range("L55").Value = 24106
$ws.Range("N55").Value = -24736
$ws.Range("H63").Value = 2181.25
$ws.Range("I63").Value = 2154.5454
$ws.Range("J63").Value = 2240
$ws.Range("K63").Value = 2154.5454
$ws.Range("L63").Value = 2240
$ws.Range("M63").Value = -1468.5454
$ws.Range("N63").Value = -3612
$ws.Range("H66").Value = 2181.25
$ws.Range("I66").Value = 2154.5454
$ws.Range("J66").Value = 2240
$ws.Range("K66").Value = 10772.727
$ws.Range("L66").Value = 11200
$ws.Range("M66").Value = -7340.726999999999
$ws.Range("N66").Value = -18064
$ws.Range("H80").Value = 34666.332
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 34666.332
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 34666.332
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -36662.332
$ws.Range("H83").Value = 34666.332
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 34666.332
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 103998.996
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -113982.996
$ws.Range("H97").Value = 510.85715
$ws.Range("I97").Value = 354.63635
$ws.Range("J97").Value = 775.2308
$ws.Range("K97").Value = 354.63635
$ws.Range("L97").Value = 775.2308
$ws.Range("M97").Value = 141.36365
$ws.Range("N97").Value = -1767.2308
$ws.Range("H102").Value = 2738.4614
$ws.Range("I102").Value = 2177.7778
$ws.Range("K102").Value = 2177.7778
$ws.Range("M102").Value = -555.7777999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 18130
$ws.Range("J35").Value = 20318.334
$ws.Range("L35").Value = 20318.334
$ws.Range("N35").Value = -20938.334
$ws.Range("H82").Value = 12254
$ws.Range("I82").Value = 1623.3334
$ws.Range("J82").Value = 28200
$ws.Range("K82").Value = 1623.3334
$ws.Range("L82").Value = 28200
$ws.Range("M82").Value = -1240.3334
$ws.Range("N82").Value = -28966
$ws.Range("H85").Value = 12254
$ws.Range("I85").Value = 1623.3334
$ws.Range("J85").Value = 28200
$ws.Range("K85").Value = 1623.3334
$ws.Range("L85").Value = 28200
$ws.Range("M85").Value = -297.3334
$ws.Range("N85").Value = -30852

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1252839.8
$ws.Range("I134").Value = 2273.2222
$ws.Range("J134").Value = 2860711
$ws.Range("K134").Value = 6819.6666
$ws.Range("L134").Value = 8582133
$ws.Range("M134").Value = -4284.6666
$ws.Range("N134").Value = -8587203

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2233572.2
$ws.Range("I5").Value = 597.6429000000001
$ws.Range("J5").Value = 4466546.5
$ws.Range("K5").Value = 1792.9287
$ws.Range("L5").Value = 13399639.5
$ws.Range("M5").Value = -1680.9287
$ws.Range("N5").Value = -13399863.5
$ws.Range("H135").Value = 2233572.2
$ws.Range("I135").Value = 597.6429000000001
$ws.Range("J135").Value = 4466546.5
$ws.Range("K135").Value = 5378.7861
$ws.Range("L135").Value = 40198918.5
$ws.Range("M135").Value = -2843.7861
$ws.Range("N135").Value = -40203988.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9105.333000000001
$ws.Range("I43").Value = 1112.5
$ws.Range("J43").Value = 18240
$ws.Range("K43").Value = 1112.5
$ws.Range("L43").Value = 18240
$ws.Range("M43").Value = -961.5
$ws.Range("N43").Value = -18542
$ws.Range("H57").Value = 19874.5
$ws.Range("I57").Value = 7500
$ws.Range("J57").Value = 23999.334
$ws.Range("K57").Value = 7500
$ws.Range("L57").Value = 23999.334
$ws.Range("M57").Value = -6680
$ws.Range("N57").Value = -25639.334
$ws.Range("H80").Value = 11979.95
$ws.Range("I80").Value = 6050
$ws.Range("J80").Value = 17909.9
$ws.Range("K80").Value = 6050
$ws.Range("L80").Value = 17909.9
$ws.Range("M80").Value = -5052
$ws.Range("N80").Value = -19905.9
$ws.Range("H83").Value = 11979.95
$ws.Range("I83").Value = 6050
$ws.Range("J83").Value = 17909.9
$ws.Range("K83").Value = 30250
$ws.Range("L83").Value = 89549.5
$ws.Range("M83").Value = -25258
$ws.Range("N83").Value = -99533.5
$ws.Range("H136").Value = 24345.688
$ws.Range("J136").Value = 24345.688
$ws.Range("L136").Value = 73037.064
$ws.Range("N136").Value = -78137.064

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5297368
$ws.Range("I132").Value = 11914570
$ws.Range("J132").Value = 3606.6
$ws.Range("K132").Value = 35743710
$ws.Range("L132").Value = 10819.8
$ws.Range("M132").Value = -35741180
$ws.Range("N132").Value = -15879.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3672150.2
$ws.Range("I132").Value = 11825
$ws.Range("K132").Value = 35475
$ws.Range("M132").Value = -32945
